$d = $word.ActiveDocument
$d.Content.Find.Execute("diciembre", $true, $false, $false, $false, $false, $true, 1, $false, "enero", 2)
Write-Host "done"
